$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, bordered, centered) from AC1 to AD1:AF1
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill season record values for each data row (2-44)
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 76
    $ws.Cells.Item($r, 31).Value = 86
    $ws.Cells.Item($r, 32).Value = 0
}
